$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing timestamp column (D2:D11) to the new extraction time
$newTimestamp = "2025-08-29 11:40:39"
for ($r = 2; $r -le 11; $r++) {
    $ws.Cells.Item($r, 4).Value = $newTimestamp
}

# New brand rows scraped from page 2
$newRows = @(
    @("MySmileUS", "189", "`$37,978,600"),
    @("American Seair Imports", "825", "`$31,078,100"),
    @("simplymandys", "82", "`$29,121,000"),
    @("Mighty Life", "12", "`$25,788,900"),
    @("Cocomint Beauty", "209", "`$25,745,000"),
    @("OQ HAIR SHOP", "118", "`$25,718,500"),
    @("SACHEU Beauty", "102", "`$25,248,400"),
    @("The Ordinary", "172", "`$24,756,600"),
    @("ONE SIZE BEAUTY", "75", "`$23,656,500"),
    @("Color Wow Hair", "56", "`$22,923,600")
)

$row = 12
foreach ($entry in $newRows) {
    $rowRange = $ws.Range("A$row`:D$row")
    # Force text storage for numeric-looking / currency-looking values so
    # Excel doesn't auto-coerce them into numbers (matches the source data,
    # which is all plain text).
    $rowRange.NumberFormat = "@"

    $ws.Cells.Item($row, 1).Value = $entry[0]
    $ws.Cells.Item($row, 2).Value = $entry[1]
    $ws.Cells.Item($row, 3).Value = $entry[2]
    $ws.Cells.Item($row, 4).Value = $newTimestamp

    # Reset the style back to the default ("Normal") so we don't leave a
    # stray number-format style applied to these cells.
    $rowRange.Style = "Normal"
    $row++
}
